$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.748.30"
$ws.Range("E2").Value = "  -4.03%  "

$ws.Range("D3").Value = "1.816.43"
$ws.Range("E3").Value = "  -3.09%  "

$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "278.66"
$ws.Range("E5").Value = "  -7.56%  "

$ws.Range("E6").Value = "  -0.17%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5090"
$ws.Range("E7").Value = "  -4.88%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3537"
$ws.Range("E8").Value = "  -5.62%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.68"
$ws.Range("E9").Value = "  -1.74%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06664"
$ws.Range("E10").Value = "  -7.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.05"
$ws.Range("E11").Value = "  -7.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.8266"
$ws.Range("E12").Value = "  -7.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07903"
$ws.Range("E13").Value = "  -3.54%  "

$ws.Range("D14").Value = "1.796.59"
$ws.Range("E14").Value = "  -4.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.075"
$ws.Range("E15").Value = "  -4.58%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.75"
$ws.Range("E16").Value = "  -5.91%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.000"
$ws.Range("E17").Value = "  -0.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.08"
$ws.Range("E18").Value = "  -5.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008034"
$ws.Range("E19").Value = "  -5.85%  "

$ws.Range("E20").Value = "  -0.15%  "

$ws.Range("D21").Value = "25.783.49"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.747"
$ws.Range("E22").Value = "  -4.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.991"
$ws.Range("E23").Value = "  -5.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.115"
$ws.Range("E24").Value = "  -4.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.233"
$ws.Range("E25").Value = "  -2.47%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.27"
$ws.Range("E26").Value = "  -2.70%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.667"
$ws.Range("E27").Value = "  -3.69%  "

$ws.Range("E28").Value = "  -5.48%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "109.33"

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.327"
$ws.Range("E30").Value = "  -8.23%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.242"
$ws.Range("E31").Value = "  -8.03%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08772"
$ws.Range("E32").Value = "  -3.81%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04908"
$ws.Range("E33").Value = "  -2.15%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7306"
$ws.Range("E34").Value = "  -9.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.139"
$ws.Range("E35").Value = "  -2.93%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.872"
$ws.Range("E36").Value = "  -2.87%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9995"
$ws.Range("E37").Value = "  -0.22%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.144"
$ws.Range("E38").Value = "  -2.27%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.393"
$ws.Range("E39").Value = "  -9.47%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01852"
$ws.Range("E40").Value = "  -5.36%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5145"
$ws.Range("E41").Value = "  -15.11%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9647"
$ws.Range("E42").Value = "  -9.99%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.219"
$ws.Range("E43").Value = "  -5.75%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "111.14"
$ws.Range("E44").Value = "  -3.32%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.035"
$ws.Range("E45").Value = "  -9.46%  "

$ws.Range("E46").Value = "  -0.15%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4560"
$ws.Range("E47").Value = "  -11.10%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1369"
$ws.Range("E48").Value = "  -8.55%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.56"
$ws.Range("E49").Value = "  -2.72%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.200"
$ws.Range("E50").Value = "  -8.02%  "

$ws.Range("E51").Value = "  -8.54%  "
